$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 81.69231000000001
$ws.Range("I5").Value = 79.625
$ws.Range("J5").Value = 106.5
$ws.Range("K5").Value = 79.625
$ws.Range("L5").Value = 106.5
$ws.Range("M5").Value = 35.375
$ws.Range("N5").Value = -336.5
$ws.Range("H18").Value = 5330.3335
$ws.Range("I18").Value = 954.1429000000001
$ws.Range("J18").Value = 11457
$ws.Range("K18").Value = 954.1429000000001
$ws.Range("L18").Value = 11457
$ws.Range("M18").Value = -670.1429000000001
$ws.Range("N18").Value = -12025
$ws.Range("H74").Value = 27786250
$ws.Range("I74").Value = 35721252
$ws.Range("K74").Value = 35721252
$ws.Range("M74").Value = -35720316
$ws.Range("H77").Value = 27786250
$ws.Range("I77").Value = 35721252
$ws.Range("K77").Value = 178606260
$ws.Range("M77").Value = -178601580

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 12837.9
$ws.Range("I31").Value = 2989.125
$ws.Range("K31").Value = 2989.125
$ws.Range("M31").Value = -2695.125
$ws.Range("H61").Value = 6797.4653
$ws.Range("I61").Value = 3112.423
$ws.Range("K61").Value = 3112.423
$ws.Range("M61").Value = -2900.423
$ws.Range("H114").Value = 59383
$ws.Range("J114").Value = 59383
$ws.Range("L114").Value = 59383
$ws.Range("N114").Value = -68061
$ws.Range("H122").Value = 2543
$ws.Range("I122").Value = 2827.4375
$ws.Range("J122").Value = 2336.1365
$ws.Range("K122").Value = 8482.3125
$ws.Range("L122").Value = 7008.4095
$ws.Range("M122").Value = -6032.3125
$ws.Range("N122").Value = -11908.4095
$ws.Range("H132").Value = 985652
$ws.Range("I132").Value = 1728090.1
$ws.Range("J132").Value = 6983.636
$ws.Range("K132").Value = 5184270.300000001
$ws.Range("L132").Value = 20950.908
$ws.Range("M132").Value = -5181740.300000001
$ws.Range("N132").Value = -26010.908
$ws.Range("H136").Value = 6797.4653
$ws.Range("I136").Value = 3112.423
$ws.Range("K136").Value = 9337.269
$ws.Range("M136").Value = -6787.269

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 11112947
$ws.Range("J20").Value = 2939.2
$ws.Range("L20").Value = 2939.2
$ws.Range("N20").Value = -3433.2
$ws.Range("H107").Value = 53576056
$ws.Range("J107").Value = 5324.6665
$ws.Range("L107").Value = 5324.6665
$ws.Range("N107").Value = -9164.666499999999
$ws.Range("H134").Value = 622.0417
$ws.Range("I134").Value = 622.0417
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 1866.1251
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 668.8749
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 364.66666
$ws.Range("I22").Value = 288.07693
$ws.Range("J22").Value = 862.5
$ws.Range("K22").Value = 288.07693
$ws.Range("L22").Value = 862.5
$ws.Range("M22").Value = 61.92307
$ws.Range("N22").Value = -1562.5
$ws.Range("H31").Value = 5745.456
$ws.Range("I31").Value = 2859.5715
$ws.Range("J31").Value = 10336.637
$ws.Range("K31").Value = 2859.5715
$ws.Range("L31").Value = 10336.637
$ws.Range("M31").Value = -2564.5715
$ws.Range("N31").Value = -10926.637
$ws.Range("H34").Value = 5745.456
$ws.Range("I34").Value = 2859.5715
$ws.Range("J34").Value = 10336.637
$ws.Range("K34").Value = 2859.5715
$ws.Range("L34").Value = 10336.637
$ws.Range("M34").Value = -2657.5715
$ws.Range("N34").Value = -10740.637
$ws.Range("H58").Value = 7170.593
$ws.Range("I58").Value = 2131.0833
$ws.Range("J58").Value = 11202.2
$ws.Range("K58").Value = 2131.0833
$ws.Range("L58").Value = 11202.2
$ws.Range("M58").Value = -1928.0833
$ws.Range("N58").Value = -11608.2
$ws.Range("H122").Value = 2383.2693
$ws.Range("I122").Value = 1650.2858
$ws.Range("J122").Value = 3238.4167
$ws.Range("K122").Value = 4950.857400000001
$ws.Range("L122").Value = 9715.250100000001
$ws.Range("M122").Value = -2500.857400000001
$ws.Range("N122").Value = -14615.2501
$ws.Range("H132").Value = 3188.6064
$ws.Range("I132").Value = 1617.9512
$ws.Range("K132").Value = 4853.8536
$ws.Range("M132").Value = -2323.8536
$ws.Range("H134").Value = 3274.877
$ws.Range("I134").Value = 1845.1177
$ws.Range("K134").Value = 5535.3531
$ws.Range("M134").Value = -3000.3531
$ws.Range("H136").Value = 7170.593
$ws.Range("I136").Value = 2131.0833
$ws.Range("J136").Value = 11202.2
$ws.Range("K136").Value = 6393.249899999999
$ws.Range("L136").Value = 33606.60000000001
$ws.Range("M136").Value = -3843.249899999999
$ws.Range("N136").Value = -38706.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100803.55
$ws.Range("I2").Value = 317.0909
$ws.Range("K2").Value = 1902.5454
$ws.Range("M2").Value = -1789.5454
$ws.Range("H56").Value = 5198545
$ws.Range("I56").Value = 5198545
$ws.Range("K56").Value = 5198545
$ws.Range("M56").Value = -5198015
$ws.Range("H68").Value = 2522.9033
$ws.Range("J68").Value = 2568.44
$ws.Range("L68").Value = 7705.32
$ws.Range("N68").Value = -9327.32
$ws.Range("H71").Value = 2522.9033
$ws.Range("J71").Value = 2568.44
$ws.Range("L71").Value = 23115.96
$ws.Range("N71").Value = -31227.96
$ws.Range("H113").Value = 2451.7646
$ws.Range("I113").Value = 1172.5
$ws.Range("J113").Value = 2845.3845
$ws.Range("K113").Value = 3517.5
$ws.Range("L113").Value = 8536.1535
$ws.Range("M113").Value = -1347.5
$ws.Range("N113").Value = -12876.1535
$ws.Range("H122").Value = 1573143.1
$ws.Range("I122").Value = 3144175.8
$ws.Range("J122").Value = 2110.4443
$ws.Range("K122").Value = 28297582.2
$ws.Range("L122").Value = 18993.9987
$ws.Range("M122").Value = -28295132.2
$ws.Range("N122").Value = -23893.9987
$ws.Range("H132").Value = 4807.5557
$ws.Range("I132").Value = 4410.8887
$ws.Range("J132").Value = 5005.8887
$ws.Range("K132").Value = 39697.99830000001
$ws.Range("L132").Value = 45052.99830000001
$ws.Range("M132").Value = -37167.99830000001
$ws.Range("N132").Value = -50112.99830000001
$ws.Range("H137").Value = 83461.96000000001
$ws.Range("J137").Value = 72205.44
$ws.Range("L137").Value = 216616.32
$ws.Range("N137").Value = -226816.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7688.5
$ws.Range("I70").Value = 5377
$ws.Range("J70").Value = 10000
$ws.Range("K70").Value = 5377
$ws.Range("L70").Value = 10000
$ws.Range("M70").Value = -5107
$ws.Range("N70").Value = -10540
$ws.Range("H73").Value = 7688.5
$ws.Range("I73").Value = 5377
$ws.Range("J73").Value = 10000
$ws.Range("K73").Value = 5377
$ws.Range("L73").Value = 10000
$ws.Range("M73").Value = -4441
$ws.Range("N73").Value = -11872
$ws.Range("H141").Value = 59990.5
$ws.Range("J141").Value = 59990.5
$ws.Range("L141").Value = 59990.5
$ws.Range("N141").Value = -70350.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 53208.332
$ws.Range("J33").Value = 54070
$ws.Range("L33").Value = 54070
$ws.Range("N33").Value = -54650
$ws.Range("H46").Value = 6176228.5
$ws.Range("J46").Value = 9262718
$ws.Range("L46").Value = 9262718
$ws.Range("N46").Value = -9263094
$ws.Range("H57").Value = 5041
$ws.Range("I57").Value = 5041
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 5041
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -4475
$ws.Range("N57").ClearContents()
$ws.Range("H100").Value = 5250.625
$ws.Range("I100").Value = 3499.5
$ws.Range("J100").Value = 7001.75
$ws.Range("K100").Value = 3499.5
$ws.Range("L100").Value = 7001.75
$ws.Range("M100").Value = -2958.5
$ws.Range("N100").Value = -8083.75
$ws.Range("H122").Value = 7760
$ws.Range("I122").Value = 8588.762000000001
$ws.Range("J122").Value = 6599.7334
$ws.Range("K122").Value = 25766.286
$ws.Range("L122").Value = 19799.2002
$ws.Range("M122").Value = -23316.286
$ws.Range("N122").Value = -24699.2002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 99988
$ws.Range("J140").Value = 99988
$ws.Range("L140").Value = 99988
$ws.Range("N140").Value = -110348
$ws.Range("H141").Value = 84985.664
$ws.Range("J141").Value = 84985.664
$ws.Range("L141").Value = 84985.664
$ws.Range("N141").Value = -95345.664

